# Edit: ppt/slides/slide1.xml - shape "CaixaDeTexto 6" (id=409)
#   - shrink the shape's height (auto-fit text box) from 321711 EMU to 290934 EMU
#   - reduce the font size of the "[local_treinamento]" paragraph's runs
#     (and its end-of-paragraph run properties) from 9pt to 7pt

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the target shape by name (falls back to the known index if needed).
$targetName = "CaixaDeTexto 6"
$shp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq $targetName) {
        $shp = $candidate
        break
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(3)
}

# --- 1. Resize the shape: <a:ext cx="4571894" cy="321711"/> -> cy="290934" ---
# 12700 EMU per point.
$shp.Height = 290934 / 12700

# --- 2. Shrink the "[local_treinamento]" run text from sz=900 (9pt) to sz=700 (7pt) ---
$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text

# The text body is:
#   Paragraph 1: "LOCAL DO TREINAMENTO: "
#   Paragraph 2: "[local_treinamento]"
# Locate paragraph 2's three runs precisely by character offsets so each
# run keeps its own rPr (bracket / placeholder name / bracket).
$openBracketPos = $fullText.IndexOf("[") + 1
$closeBracketPos = $fullText.IndexOf("]") + 1
$nameStart = $openBracketPos + 1
$nameLength = $closeBracketPos - $nameStart

$runOpen = $tr.Characters($openBracketPos, 1)
$runOpen.Font.Size = 7

$runName = $tr.Characters($nameStart, $nameLength)
$runName.Font.Size = 7

$runClose = $tr.Characters($closeBracketPos, 1)
$runClose.Font.Size = 7
